# Regenerate orders with updated distance/sizes:
#   Distance codes: D64 -> D69, D51 -> D55, D80 -> D86
#   Size code:      S30 -> S31
# These codes appear both standalone (Distance / Size columns) and embedded
# inside composite strings (Condition, Filename_Left, Filename_Right), so we
# do a partial (substring) find & replace across the whole used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

$xlPart = 2   # xlPart: match substrings, not just whole-cell contents

$used.Replace("D64", "D69", $xlPart)
$used.Replace("D51", "D55", $xlPart)
$used.Replace("D80", "D86", $xlPart)
$used.Replace("S30", "S31", $xlPart)
